# Updates room allocations ("sala") in column H of the "Alocações" sheet
# so that each room is used by only one allocation at a time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Alocações")

# Map of row number -> new value for column H
$updates = @{
    2  = "['210-B']"
    3  = "['210-B']"
    4  = "['210-B']"
    5  = "['210-B']"
    6  = "['301-A']"
    7  = "['301-A']"
    8  = "['301-A']"
    9  = "['302-A']"
    10 = "['302-A']"
    11 = "['301-A']"
    12 = "['301-A']"
    13 = "['301-A']"
    14 = "['301-A']"
    15 = "['301-A']"
    16 = "['301-A']"
    17 = "['301-A']"
    18 = "['305-B']"
    19 = "['304-B']"
    20 = "['304-B']"
    21 = "['304-B']"
    22 = "['305-B']"
    23 = "['304-B']"
    25 = "['302-B']"
    26 = "['303-B']"
    27 = "['302-B']"
    28 = "['302-B']"
    29 = "['302-B']"
    30 = "['303-B']"
    32 = "['309-A']"
    33 = "['304-A']"
    34 = "['305-A']"
    35 = "['305-A']"
    36 = "['305-A']"
    37 = "['305-A']"
    38 = "['206-B']"
    39 = "['206-B']"
    40 = "['206-B']"
    41 = "['206-B']"
    42 = "['206-B']"
    43 = "['206-B']"
    44 = "['206-B']"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 8).Value = $updates[$row]
}

$wb.Save()
